# Motor_pos.xlsx - Simulink step-input lookup table was shifted/edited:
# a block of rows had their B (torque?) and C (speed?) step values updated
# to re-time the waveform used for the SIMX/Simulink co-simulation test.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1-11: step starts at -10 / 30 instead of 0 / 20
$ws.Range("B1:B11").Value = -10
$ws.Range("C1:C11").Value = 30

# Rows 32-41: C drops from 40 to 30 (the -10/40 plateau now starts one block later)
$ws.Range("C32:C41").Value = 30
# Rows 72-81: C rises from 30 to 40 (the -10/30 plateau now starts one block later)
$ws.Range("C72:C81").Value = 40

# Rows 112-121: B drops from 0 to -10 (the 0-torque plateau now starts one block later)
$ws.Range("B112:B121").Value = -10

# Rows 132-141: C drops from 40 to 30
$ws.Range("C132:C141").Value = 30
# Rows 172-181: C rises from 30 to 40
$ws.Range("C172:C181").Value = 40

# Rows 212-221: B drops from 10 to 0 (the 10-torque plateau now starts one block later)
$ws.Range("B212:B221").Value = 0

# Rows 232-241: C drops from 40 to 30
$ws.Range("C232:C241").Value = 30
# Rows 272-281: C rises from 30 to 40
$ws.Range("C272:C281").Value = 40

# Rows 312-331 (tail of the sheet): B rises from 0 to 10, C rises from 20 to 30
$ws.Range("B312:B331").Value = 10
$ws.Range("C312:C331").Value = 30

# Reflect the author's final scroll/selection position in the sheet view
$ws.Activate()
$ws.Range("A322:D331").Select()
$excel.ActiveWindow.ScrollRow = 316
$excel.ActiveWindow.ScrollColumn = 1
